$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 25.6
$ws.Range("I6").Value = 7
$ws.Range("K6").Value = 21
$ws.Range("M6").Value = 91
$ws.Range("H15").Value = 2308.875
$ws.Range("I15").Value = 2308.875
$ws.Range("K15").Value = 6926.625
$ws.Range("M15").Value = -6757.625
$ws.Range("H62").Value = 8906.4
$ws.Range("I62").Value = 8783
$ws.Range("K62").Value = 8783
$ws.Range("M62").Value = -8159
$ws.Range("H65").Value = 8906.4
$ws.Range("I65").Value = 8783
$ws.Range("K65").Value = 43915
$ws.Range("M65").Value = -40795
$ws.Range("H69").Value = 375
$ws.Range("I69").Value = 375
$ws.Range("K69").Value = 1125
$ws.Range("M69").Value = -251
$ws.Range("H72").Value = 375
$ws.Range("I72").Value = 375
$ws.Range("K72").Value = 3375
$ws.Range("M72").Value = 993
$ws.Range("H96").Value = 1045
$ws.Range("I96").Value = 1045
$ws.Range("K96").Value = 3135
$ws.Range("M96").Value = -1762
$ws.Range("H129").Value = 1185.4445
$ws.Range("I129").Value = 1185.4445
$ws.Range("K129").Value = 3556.3335
$ws.Range("M129").Value = 1443.6665
$ws.Range("H141").Value = 24640
$ws.Range("I141").Value = 15568
$ws.Range("K141").Value = 46704
$ws.Range("M141").Value = -41524
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3070.5
$ws.Range("I32").Value = 2883.276
$ws.Range("K32").Value = 2883.276
$ws.Range("M32").Value = -2596.276
$ws.Range("H61").Value = 4000
$ws.Range("J61").Value = 4000
$ws.Range("L61").Value = 4000
$ws.Range("N61").Value = -4424
$ws.Range("H74").Value = 4274.5
$ws.Range("I74").Value = 4274.5
$ws.Range("K74").Value = 4274.5
$ws.Range("M74").Value = -3400.5
$ws.Range("H77").Value = 4274.5
$ws.Range("I77").Value = 4274.5
$ws.Range("K77").Value = 21372.5
$ws.Range("M77").Value = -17004.5
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H97").Value = 1804.1818
$ws.Range("I97").Value = 855.875
$ws.Range("K97").Value = 855.875
$ws.Range("M97").Value = -359.875
$ws.Range("H132").Value = 2399.5
$ws.Range("I132").Value = 2399.5
$ws.Range("K132").Value = 7198.5
$ws.Range("M132").Value = -4668.5
$ws.Range("H136").Value = 4000
$ws.Range("J136").Value = 4000
$ws.Range("L136").Value = 12000
$ws.Range("N136").Value = -17100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 20000
$ws.Range("I99").Value = 20000
$ws.Range("K99").Value = 20000
$ws.Range("M99").Value = -18502
$ws.Range("H105").Value = 3930.5
$ws.Range("I105").Value = 3930.5
$ws.Range("K105").Value = 3930.5
$ws.Range("M105").Value = -2183.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 841.6667
$ws.Range("I35").Value = 841.6667
$ws.Range("K35").Value = 841.6667
$ws.Range("M35").Value = -547.6667
$ws.Range("H132").Value = 7266.3335
$ws.Range("I132").Value = 6999.5
$ws.Range("K132").Value = 20998.5
$ws.Range("M132").Value = -18468.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 999999
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 999999
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 2999997
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -3000221
$ws.Range("H33").Value = 159.75
$ws.Range("I33").Value = 24.5
$ws.Range("K33").Value = 147
$ws.Range("M33").Value = 136
$ws.Range("H92").Value = 5000
$ws.Range("I92").Value = 5000
$ws.Range("J92").Value = 5000
$ws.Range("K92").Value = 15000
$ws.Range("L92").Value = 15000
$ws.Range("M92").Value = -13752
$ws.Range("N92").Value = -17496
$ws.Range("H94").Value = 2549
$ws.Range("I94").Value = 300
$ws.Range("J94").Value = 2998.8
$ws.Range("K94").Value = 900
$ws.Range("L94").Value = 8996.400000000001
$ws.Range("M94").Value = -224
$ws.Range("N94").Value = -10348.4
$ws.Range("H97").Value = 400
$ws.Range("I97").Value = 400
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1200
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -704
$ws.Range("N97").ClearContents()
$ws.Range("H99").Value = 2950
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 885.2857
$ws.Range("I107").Value = 900
$ws.Range("J107").Value = 882.8333
$ws.Range("K107").Value = 2700
$ws.Range("L107").Value = 2648.4999
$ws.Range("M107").Value = -780
$ws.Range("N107").Value = -6488.4999
$ws.Range("H110").Value = 499.5
$ws.Range("I110").Value = 499.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1498.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 2591.5
$ws.Range("N110").ClearContents()
$ws.Range("H112").Value = 857.5
$ws.Range("I112").Value = 857.5
$ws.Range("K112").Value = 2572.5
$ws.Range("M112").Value = -1464.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 7000399.5
$ws.Range("I14").Value = 7500250
$ws.Range("J14").Value = 6667166
$ws.Range("K14").Value = 7500250
$ws.Range("L14").Value = 6667166
$ws.Range("M14").Value = -7500082
$ws.Range("N14").Value = -6667502
$ws.Range("H132").Value = 4370.1113
$ws.Range("I132").Value = 4291.375
$ws.Range("K132").Value = 12874.125
$ws.Range("M132").Value = -10344.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1733.4166
$ws.Range("J46").Value = 2091.125
$ws.Range("L46").Value = 2091.125
$ws.Range("N46").Value = -2467.125
$ws.Range("H100").Value = 4499.6
$ws.Range("I100").Value = 4999
$ws.Range("J100").Value = 4166.6665
$ws.Range("K100").Value = 4999
$ws.Range("L100").Value = 4166.6665
$ws.Range("M100").Value = -4458
$ws.Range("N100").Value = -5248.6665
$ws.Range("H122").Value = 3692.1667
$ws.Range("I122").Value = 3538.75
$ws.Range("K122").Value = 10616.25
$ws.Range("M122").Value = -8166.25
$ws.Range("H130").Value = 8284.200000000001
$ws.Range("J130").Value = 8284.200000000001
$ws.Range("L130").Value = 8284.200000000001
$ws.Range("N130").Value = -18324.2
$ws.Range("H132").Value = 19707.77
$ws.Range("I132").Value = 20644
$ws.Range("K132").Value = 61932
$ws.Range("M132").Value = -59402
$ws.Range("H136").Value = 4438
$ws.Range("I136").Value = 4072
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 12216
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -9666
$ws.Range("N136").Value = -26100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 30003
$ws.Range("I3").Value = 30003
$ws.Range("K3").Value = 30003
$ws.Range("M3").Value = -29889
$ws.Range("H132").Value = 3523
$ws.Range("I132").Value = 3523
$ws.Range("K132").Value = 10569
$ws.Range("M132").Value = -8039
$ws.Range("H136").Value = 3483
$ws.Range("I136").Value = 2998.25
$ws.Range("J136").Value = 4452.5
$ws.Range("K136").Value = 8994.75
$ws.Range("L136").Value = 13357.5
$ws.Range("M136").Value = -6444.75
$ws.Range("N136").Value = -18457.5
